# dsAlgoData.xlsx - update the SignIn test-data credentials used by the
# NinjaTester Excel Reader (ConfigReader) and add the hyperlink Excel
# applies automatically to the new "Welcome@123" value.

$wb = $excel.ActiveWorkbook
$signIn = $wb.Worksheets.Item("SignIn")

# New username / password test data (was "kodetesters" / "numpyninja24").
$signIn.Range("A2").Value = "NinjaTester"
$signIn.Range("B2").Value = "Welcome@123"

# Excel auto-links values containing "@"; recreate that hyperlink on B2.
$signIn.Hyperlinks.Add($signIn.Range("B2"), "mailto:Welcome@123")
